$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1195.8334
$ws.Range("J32").Value = 1625
$ws.Range("L32").Value = 1625
$ws.Range("N32").Value = -2277

$ws.Range("H64").Value = 4293.3335
$ws.Range("I64").Value = 4495
$ws.Range("J64").Value = 3890
$ws.Range("K64").Value = 4495
$ws.Range("L64").Value = 3890
$ws.Range("M64").Value = -4247
$ws.Range("N64").Value = -4386

$ws.Range("H67").Value = 4293.3335
$ws.Range("I67").Value = 4495
$ws.Range("J67").Value = 3890
$ws.Range("K67").Value = 4495
$ws.Range("L67").Value = 3890
$ws.Range("M67").Value = -3637
$ws.Range("N67").Value = -5606

$ws.Range("H98").Value = 10271.429
$ws.Range("I98").Value = 11816.667
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 11816.667
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = -10318.667
$ws.Range("N98").Value = -3996

$ws.Range("H113").Value = 1125
$ws.Range("I113").Value = 1500
$ws.Range("K113").Value = 1500
$ws.Range("M113").Value = 1754

$ws.Range("H122").Value = 10271.429
$ws.Range("I122").Value = 11816.667
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 35450.001
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -33000.001
$ws.Range("N122").Value = -7900

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1879.875
$ws.Range("I110").Value = 1475
$ws.Range("J110").Value = 2014.8334
$ws.Range("K110").Value = 1475
$ws.Range("L110").Value = 2014.8334
$ws.Range("M110").Value = 570
$ws.Range("N110").Value = -6104.8334

$ws.Range("H122").Value = 2103.4375
$ws.Range("I122").Value = 1685.5834
$ws.Range("K122").Value = 5056.7502
$ws.Range("M122").Value = -2606.7502

$ws.Range("H132").Value = 1934.9272
$ws.Range("I132").Value = 1536
$ws.Range("K132").Value = 4608
$ws.Range("M132").Value = -2078

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 19231654
$ws.Range("I94").Value = 35714844
$ws.Range("J94").Value = 1265
$ws.Range("K94").Value = 35714844
$ws.Range("L94").Value = 1265
$ws.Range("M94").Value = -35714393
$ws.Range("N94").Value = -2167

$ws.Range("H105").Value = 142858800
$ws.Range("I105").Value = 166668240
$ws.Range("K105").Value = 166668240
$ws.Range("M105").Value = -166666493

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 66668064
$ws.Range("I16").Value = 90910380
$ws.Range("K16").Value = 90910380
$ws.Range("M16").Value = -90910093

$ws.Range("H22").Value = 140440.4
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 233734
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 233734
$ws.Range("M22").Value = -150
$ws.Range("N22").Value = -234434

$ws.Range("H31").Value = 1765.1702
$ws.Range("I31").Value = 1439.6904
$ws.Range("J31").Value = 4499.2
$ws.Range("K31").Value = 1439.6904
$ws.Range("L31").Value = 4499.2
$ws.Range("M31").Value = -1144.6904
$ws.Range("N31").Value = -5089.2

$ws.Range("H34").Value = 1765.1702
$ws.Range("I34").Value = 1439.6904
$ws.Range("J34").Value = 4499.2
$ws.Range("K34").Value = 1439.6904
$ws.Range("L34").Value = 4499.2
$ws.Range("M34").Value = -1237.6904
$ws.Range("N34").Value = -4903.2

$ws.Range("H113").Value = 66668064
$ws.Range("I113").Value = 90910380
$ws.Range("K113").Value = 90910380
$ws.Range("M113").Value = -90908210

$ws.Range("H122").Value = 788.55
$ws.Range("I122").Value = 731.7222
$ws.Range("K122").Value = 2195.1666
$ws.Range("M122").Value = 254.8334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 378.85715
$ws.Range("I86").Value = 152.85715
$ws.Range("J86").Value = 604.8570999999999
$ws.Range("K86").Value = 458.57145
$ws.Range("L86").Value = 1814.5713
$ws.Range("M86").Value = 727.4285500000001
$ws.Range("N86").Value = -4186.5713

$ws.Range("H89").Value = 378.85715
$ws.Range("I89").Value = 152.85715
$ws.Range("J89").Value = 604.8570999999999
$ws.Range("K89").Value = 1375.71435
$ws.Range("L89").Value = 5443.7139
$ws.Range("M89").Value = 4552.28565
$ws.Range("N89").Value = -17299.7139

$ws.Range("H139").Value = 2267.6177
$ws.Range("I139").Value = 2925.2
$ws.Range("J139").Value = 1748.4736
$ws.Range("K139").Value = 8775.599999999999
$ws.Range("L139").Value = 5245.4208
$ws.Range("M139").Value = -3635.599999999999
$ws.Range("N139").Value = -15525.4208

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1446.7646
$ws.Range("I102").Value = 1171.0834
$ws.Range("J102").Value = 2108.4
$ws.Range("K102").Value = 1171.0834
$ws.Range("L102").Value = 2108.4
$ws.Range("M102").Value = 450.9166
$ws.Range("N102").Value = -5352.4

$ws.Range("H122").Value = 3562.4
$ws.Range("I122").Value = 3754.353
$ws.Range("J122").Value = 3154.5
$ws.Range("K122").Value = 11263.059
$ws.Range("L122").Value = 9463.5
$ws.Range("M122").Value = -8813.059000000001
$ws.Range("N122").Value = -14363.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3764.2273
$ws.Range("I40").Value = 1904
$ws.Range("K40").Value = 1904
$ws.Range("M40").Value = -1768

$ws.Range("H61").Value = 1154.1111
$ws.Range("I61").Value = 1097
$ws.Range("J61").Value = 1268.3334
$ws.Range("K61").Value = 1097
$ws.Range("L61").Value = 1268.3334
$ws.Range("M61").Value = -895
$ws.Range("N61").Value = -1672.3334

$ws.Range("H113").Value = 1154.1111
$ws.Range("I113").Value = 1097
$ws.Range("J113").Value = 1268.3334
$ws.Range("K113").Value = 1097
$ws.Range("L113").Value = 1268.3334
$ws.Range("M113").Value = 1073
$ws.Range("N113").Value = -5608.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2192.55
$ws.Range("I96").Value = 1988.4615
$ws.Range("J96").Value = 2571.5715
$ws.Range("K96").Value = 1988.4615
$ws.Range("L96").Value = 2571.5715
$ws.Range("M96").Value = -615.4614999999999
$ws.Range("N96").Value = -5317.5715

$ws.Range("H122").Value = 10417937
$ws.Range("I122").Value = 11906090
$ws.Range("J122").Value = 868.3333
$ws.Range("K122").Value = 35718270
$ws.Range("L122").Value = 2604.9999
$ws.Range("M122").Value = -35715820
$ws.Range("N122").Value = -7504.9999

$ws.Range("H126").Value = 58823830
$ws.Range("I126").Value = 58823830
$ws.Range("K126").Value = 176471490
$ws.Range("M126").Value = -176469020
